$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.694.09"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.632.88"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.90"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.497"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0619"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.96"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "1.860.84"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.617.20"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.81"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "26.665.32"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  -3.27%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.26"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.29"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.33"
$ws.Range("E22").Value = "  -8.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.15"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.32"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.02"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").Value = "1.259.68"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0172"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.519"
$ws.Range("E38").Value = "  -3.98%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.796"
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").Value = "1.772.85"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.23"
$ws.Range("E44").Value = "  -3.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.85"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.45"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0951"
$ws.Range("E51").Value = "  -2.70%  "
